$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @(D_new, E_new)  ; D_new = $null means column D unchanged
$changes = @{
    2  = @("61.476.36", "  -2.28%  ")
    3  = @("3.380.45",  "  -0.80%  ")
    4  = @("1.00",      "  -0.03%  ")
    5  = @("577.69",    "  +0.27%  ")
    6  = @("136.20",    "  +7.35%  ")
    7  = @($null,       "  -0.04%  ")
    8  = @("3.380.40",  "  -0.82%  ")
    9  = @("0.479",     "  +0.29%  ")
    10 = @("7.59",      "  +2.47%  ")
    11 = @("0.124",     "  +1.74%  ")
    12 = @("0.389",     "  +1.96%  ")
    13 = @("3.955.91",  "  -0.89%  ")
    14 = @($null,       "  +1.13%  ")
    15 = @("0.0000176", "  +1.11%  ")
    16 = @("3.377.18",  "  -0.89%  ")
    17 = @("25.35",     "  +1.24%  ")
    18 = @("61.543.62", "  -2.28%  ")
    19 = @("14.09",     "  +6.58%  ")
    20 = @("5.83",      "  +1.94%  ")
    21 = @("9.37",      "  -2.35%  ")
    22 = @("385.43",    "  +1.90%  ")
    23 = @("0.569",     "  +1.61%  ")
    24 = @("3.513.40",  "  -0.91%  ")
    25 = @($null,       "  +0.16%  ")
    26 = @("70.97",     "  -2.23%  ")
    27 = @("0.0000120", "  +10.01%  ")
    28 = @("1.71",      "  +20.82%  ")
    29 = @("7.86",      "  +12.03%  ")
    30 = @("0.996",     "  -0.50%  ")
    31 = @("8.18",      "  +3.32%  ")
    32 = @($null,       "  +0.59%  ")
    33 = @("0.157",     "  +2.93%  ")
    34 = @($null,       "  -0.04%  ")
    35 = @("3.412.63",  "  -0.66%  ")
    36 = @("23.46",     "  +2.48%  ")
    37 = @("5.61",      "  +5.21%  ")
    38 = @("7.01",      "  +3.41%  ")
    39 = @($null,       "  +3.33%  ")
    40 = @("163.02",    "  -0.57%  ")
    41 = @("0.0787",    "  +3.07%  ")
    42 = @($null,       "  +13.52%  ")
    43 = @($null,       "  +0.01%  ")
    44 = @($null,       "  +3.25%  ")
    45 = @("41.77",     "  +0.28%  ")
    46 = @("0.759",     "  -2.64%  ")
    47 = @($null,       "  +2.29%  ")
    48 = @("23.55",     "  +2.77%  ")
    49 = @($null,       "  +3.78%  ")
    50 = @("23.44",     "  +15.31%  ")
    51 = @("0.905",     "  +5.02%  ")
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $dNew = $vals[0]
    $eNew = $vals[1]

    if ($null -ne $dNew) {
        # Column D values are strings such as "61.476.36" or "1.00" that Excel
        # would otherwise auto-parse as numbers/dates. Force text entry by
        # temporarily formatting as Text, then strip the formatting again so
        # the cell ends up with no explicit style, same as the source file.
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $dNew
        $dCell.ClearFormats()
    }

    # Column E values (e.g. "  -2.28%  ") already contain spaces/percent signs
    # so Excel keeps them as plain text; Value2 avoids any locale reparsing.
    $ws.Cells.Item($row, 5).Value2 = $eNew
}
